# Urban population by ethnicity
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header renames / new headers -------------------------------------
# Order matters: it controls the position new strings land in the shared
# string table, matching the original authoring order.
$ws.Range("F1").Value = "urb_rate"
$ws.Range("D1").Value = "urban_pop"
$ws.Range("B1").Value = "pop_total"
$ws.Range("E1").Value = "rural_pop"

# --- New urban population column (D) -----------------------------------
$ws.Cells.Item(2, 4).Value = 555675
$ws.Cells.Item(3, 4).Value = 489115
$ws.Cells.Item(4, 4).Value = 35578
$ws.Cells.Item(5, 4).Value = 184769
$ws.Cells.Item(6, 4).Value = 45347
$ws.Cells.Item(7, 4).Value = 3111
$ws.Cells.Item(8, 4).Value = 6647
$ws.Cells.Item(9, 4).Value = 2143262
$ws.Cells.Item(10, 4).Value = 1781
$ws.Cells.Item(11, 4).Value = 1670
$ws.Cells.Item(12, 4).Value = 9156
$ws.Cells.Item(13, 4).Value = 59959
$ws.Cells.Item(14, 4).Value = 26190
$ws.Cells.Item(15, 4).Value = 21392
$ws.Cells.Item(16, 4).Value = 256185
$ws.Cells.Item(17, 4).Value = 451440
$ws.Cells.Item(18, 4).Value = 3286945

# --- Rural population column (E) = pop_total - urban_pop ----------------
$ws.Range("E2:E15").Formula = "=B2-D2"
# --- Urbanization rate column (F) = urban_pop / pop_total ----------------
$ws.Range("F2").Formula = "=D2/B2"
$ws.Range("F3:F18").Formula = "=D3/B3"
$ws.Range("E16").Formula = "=B16-D16"
$ws.Range("E17:E18").Formula = "=B17-D17"

# Referencing B4/B18 (which carry a custom number format) makes the engine
# auto-inherit that number format onto the new formula cells; the source
# workbook does not show that, so strip the picked-up formatting back off.
$ws.Range("E4").ClearFormats()
$ws.Range("E18").ClearFormats()

# --- Font styling for the pasted D18 value (Wikipedia-style paste) -------
$d18 = $ws.Range("D18")
$d18.Font.Name = "Arial"
$d18.Font.Size = 8
$d18.Font.Color = 1710618

# --- Percent-style formatting for the rate column ------------------------
$ws.Range("F2:F18").Style = "Percent"

# --- Selection / view state ------------------------------------------
$ws.Range("F5").Select() | Out-Null
